$d = $word.ActiveDocument

$replacements = @(
    @{old="84÷7="; new="20÷5="},
    @{old="95÷9="; new="73÷5="},
    @{old="25÷9="; new="87÷8="},
    @{old="63÷9="; new="36÷3="},
    @{old="98÷3="; new="90÷5="},
    @{old="35÷2="; new="57÷8="},
    @{old="29÷9="; new="75÷9="},
    @{old="31÷2="; new="37÷7="},
    @{old="17÷6="; new="76÷5="},
    @{old="36÷4="; new="97÷9="},
    @{old="34÷8="; new="64÷9="},
    @{old="82÷9="; new="73÷9="},
    @{old="52÷9="; new="96÷4="},
    @{old="36÷6="; new="61÷9="},
    @{old="66÷7="; new="72÷9="},
    @{old="71÷7="; new="68÷9="},
    @{old="87÷5="; new="30÷5="},
    @{old="68÷2="; new="66÷5="},
    @{old="68÷5="; new="15÷7="},
    @{old="74÷4="; new="18÷7="},
    @{old="53÷8="; new="28÷9="},
    @{old="56÷2="; new="39÷6="},
    @{old="80÷8="; new="43÷6="},
    @{old="75÷8="; new="15÷2="},
    @{old="71÷6="; new="84÷9="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
